$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Code Review 1 scores for each contributor
$ws.Range("B2").Value = 20
$ws.Range("B3").Value = 20
$ws.Range("B4").Value = 20
$ws.Range("B5").Value = 20
$ws.Range("B6").Value = 20

# Update the active selection to match the saved cursor position
$ws.Range("C11").Select()
